# Apply weekly update to the Brocoli / Macroferia Regional de Talca sheet.
# Net effect (per the xml diff): two brand-new data rows are inserted into
# the table, pushing the existing rows below them down by one each time.
#   - a new row is inserted at row 314 (pushes old 314..351 -> 315..352)
#   - a new row is inserted at (the then-current) row 342
#     (pushes old 341(now at 342)..351(now at 352) -> 343..353)
# Dimension grows from A1:R351 to A1:R353.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row (final row 314) ---
$ws.Rows.Item(314).Insert()

$ws.Cells.Item(314, 1).Value = 5
$ws.Cells.Item(314, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(314, 3).Value = "Maule"
$ws.Cells.Item(314, 4).Value = 44748
$ws.Cells.Item(314, 5).Value = 7
$ws.Cells.Item(314, 6).Value = 100112023
$ws.Cells.Item(314, 7).Value = "Brócoli"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 5000
$ws.Cells.Item(314, 11).Value = 700
$ws.Cells.Item(314, 12).Value = 700
$ws.Cells.Item(314, 13).Value = 700
$ws.Cells.Item(314, 14).Value = "$/unidad"
$ws.Cells.Item(314, 15).Value = "Región del Maule"
$ws.Cells.Item(314, 16).Value = 700
$ws.Cells.Item(314, 17).Value = 1
$ws.Cells.Item(314, 18).Value = "Hortaliza"

# --- Insert second new row (final row 342) ---
$ws.Rows.Item(342).Insert()

$ws.Cells.Item(342, 1).Value = 5
$ws.Cells.Item(342, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(342, 3).Value = "Maule"
$ws.Cells.Item(342, 4).Value = 44747
$ws.Cells.Item(342, 5).Value = 7
$ws.Cells.Item(342, 6).Value = 100112023
$ws.Cells.Item(342, 7).Value = "Brócoli"
$ws.Cells.Item(342, 8).Value = "Sin especificar"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 5000
$ws.Cells.Item(342, 11).Value = 800
$ws.Cells.Item(342, 12).Value = 800
$ws.Cells.Item(342, 13).Value = 800
$ws.Cells.Item(342, 14).Value = "$/unidad"
$ws.Cells.Item(342, 15).Value = "Región del Maule"
$ws.Cells.Item(342, 16).Value = 800
$ws.Cells.Item(342, 17).Value = 1
$ws.Cells.Item(342, 18).Value = "Hortaliza"

# --- Make sure column D keeps the date-ish number format used by the rest
#     of the "Fecha" column for the two freshly-inserted rows. ---
$ws.Range("D314").NumberFormat = $ws.Range("D313").NumberFormat
$ws.Range("D342").NumberFormat = $ws.Range("D313").NumberFormat
